$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 14.99
$ws.Range("C5").Value = 1.01
$ws.Range("D10").Value = 4586

$ws.Range("A1:D13").Borders.LineStyle = 1
